$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$csvData = @"
2,45969,0.37,0,1,08.11.20251
3,45969.01041666666,0,0,2,08.11.20252
4,45969.02083333334,0,0,3,08.11.20253
5,45969.03125,0,0,4,08.11.20254
6,45969.04166666666,0.45,0,5,08.11.20255
7,45969.05208333334,0,0,6,08.11.20256
8,45969.0625,0,0,7,08.11.20257
9,45969.07291666666,0,0,8,08.11.20258
10,45969.08333333334,0,0,9,08.11.20259
11,45969.09375,0,0,10,08.11.202510
12,45969.10416666666,0,0,11,08.11.202511
13,45969.11458333334,0,0,12,08.11.202512
14,45969.125,0.53,0,13,08.11.202513
15,45969.13541666666,0,0,14,08.11.202514
16,45969.14583333334,0,0,15,08.11.202515
17,45969.15625,0,0,16,08.11.202516
18,45969.16666666666,0.45,0,17,08.11.202517
19,45969.17708333334,0,0,18,08.11.202518
20,45969.1875,0,0,19,08.11.202519
21,45969.19791666666,0,0,20,08.11.202520
22,45969.20833333334,0.423,0,21,08.11.202521
23,45969.21875,0.495,0,22,08.11.202522
24,45969.22916666666,0.607,0,23,08.11.202523
25,45969.23958333334,0.763,0,24,08.11.202524
26,45969.25,19.75,0,25,08.11.202525
27,45969.26041666666,25.674,0,26,08.11.202526
28,45969.27083333334,34.462,9,27,08.11.202527
29,45969.28125,48.081,34,28,08.11.202528
30,45969.29166666666,130.202,70,29,08.11.202529
31,45969.30208333334,154.602,0,30,08.11.202530
32,45969.3125,180.326,0,31,08.11.202531
33,45969.32291666666,203.857,0,32,08.11.202532
34,45969.33333333334,356.773,0,33,08.11.202533
35,45969.34375,385.956,0,34,08.11.202534
36,45969.35416666666,429.691,0,35,08.11.202535
37,45969.36458333334,454.672,0,36,08.11.202536
38,45969.375,504.809,0,37,08.11.202537
39,45969.38541666666,527.515,0,38,08.11.202538
40,45969.39583333334,546.627,0,39,08.11.202539
41,45969.40625,564.298,0,40,08.11.202540
42,45969.41666666666,568.2910000000001,0,41,08.11.202541
43,45969.42708333334,579.2670000000001,0,42,08.11.202542
44,45969.4375,587.7089999999999,0,43,08.11.202543
45,45969.44791666666,589.253,0,44,08.11.202544
46,45969.45833333334,572.551,0,45,08.11.202545
47,45969.46875,566.65,0,46,08.11.202546
48,45969.47916666666,555.53,0,47,08.11.202547
49,45969.48958333334,537.359,0,48,08.11.202548
50,45969.5,478.003,0,49,08.11.202549
51,45969.51041666666,457.065,0,50,08.11.202550
52,45969.52083333334,437.962,0,51,08.11.202551
53,45969.53125,416.554,0,52,08.11.202552
54,45969.54166666666,349.649,0,53,08.11.202553
55,45969.55208333334,324.558,0,54,08.11.202554
56,45969.5625,297.557,0,55,08.11.202555
57,45969.57291666666,272.195,0,56,08.11.202556
58,45969.58333333334,187.41,0,57,08.11.202557
59,45969.59375,162.861,0,58,08.11.202558
60,45969.60416666666,135.766,0,59,08.11.202559
61,45969.61458333334,116.332,0,60,08.11.202560
62,45969.625,61.702,0,61,08.11.202561
63,45969.63541666666,46.725,0,62,08.11.202562
64,45969.64583333334,34.73,0,63,08.11.202563
65,45969.65625,30.651,0,64,08.11.202564
66,45969.66666666666,18.062,0,65,08.11.202565
67,45969.67708333334,17.315,0,66,08.11.202566
68,45969.6875,14.958,0,67,08.11.202567
69,45969.69791666666,14.509,0,68,08.11.202568
70,45969.70833333334,1.29,0,69,08.11.202569
71,45969.71875,1.158,0,70,08.11.202570
72,45969.72916666666,0.954,0,71,08.11.202571
73,45969.73958333334,0.914,0,72,08.11.202572
74,45969.75,0.386,0,73,08.11.202573
75,45969.76041666666,0,0,74,08.11.202574
76,45969.77083333334,0,0,75,08.11.202575
77,45969.78125,0,0,76,08.11.202576
78,45969.79166666666,0.45,0,77,08.11.202577
79,45969.80208333334,0,0,78,08.11.202578
80,45969.8125,0,0,79,08.11.202579
81,45969.82291666666,0,0,80,08.11.202580
82,45969.83333333334,0,0,81,08.11.202581
83,45969.84375,0,0,82,08.11.202582
84,45969.85416666666,0,0,83,08.11.202583
85,45969.86458333334,0,0,84,08.11.202584
86,45969.875,0.53,0,85,08.11.202585
87,45969.88541666666,0,0,86,08.11.202586
88,45969.89583333334,0,0,87,08.11.202587
89,45969.90625,0,0,88,08.11.202588
90,45969.91666666666,0.45,0,89,08.11.202589
91,45969.92708333334,0,0,90,08.11.202590
92,45969.9375,0,0,91,08.11.202591
93,45969.94791666666,0,0,92,08.11.202592
94,45969.95833333334,0,0,93,08.11.202593
95,45969.96875,0,0,94,08.11.202594
96,45969.97916666666,0,0,95,08.11.202595
97,45969.98958333334,0,0,96,08.11.202596
98,45970,0.58,0.0,1,09.11.20251
99,45970.01041666666,0.0,0.0,2,09.11.20252
100,45970.02083333334,0.0,0.0,3,09.11.20253
101,45970.03125,0.0,0.0,4,09.11.20254
102,45970.04166666666,0.54,0.0,5,09.11.20255
103,45970.05208333334,0.0,0.0,6,09.11.20256
104,45970.0625,0.0,0.0,7,09.11.20257
105,45970.07291666666,0.0,0.0,8,09.11.20258
106,45970.08333333334,0.0,0.0,9,09.11.20259
107,45970.09375,0.0,0.0,10,09.11.202510
108,45970.10416666666,0.0,0.0,11,09.11.202511
109,45970.11458333334,0.0,0.0,12,09.11.202512
110,45970.125,0.0,0.0,13,09.11.202513
111,45970.13541666666,0.0,0.0,14,09.11.202514
112,45970.14583333334,0.0,0.0,15,09.11.202515
113,45970.15625,0.0,0.0,16,09.11.202516
114,45970.16666666666,0.58,0.0,17,09.11.202517
115,45970.17708333334,0.0,0.0,18,09.11.202518
116,45970.1875,0.0,0.0,19,09.11.202519
117,45970.19791666666,0.0,0.0,20,09.11.202520
118,45970.20833333334,0.719,0.0,21,09.11.202521
119,45970.21875,0.817,0.0,22,09.11.202522
120,45970.22916666666,1.144,0.0,23,09.11.202523
121,45970.23958333334,1.592,0.0,24,09.11.202524
122,45970.25,15.346,0.0,25,09.11.202525
123,45970.26041666666,22.35,0.0,26,09.11.202526
124,45970.27083333334,34.385,7.0,27,09.11.202527
125,45970.28125,49.305,21.0,28,09.11.202528
126,45970.29166666666,149.017,51.0,29,09.11.202529
127,45970.30208333334,178.118,0.0,30,09.11.202530
128,45970.3125,209.611,0.0,31,09.11.202531
129,45970.32291666666,244.246,0.0,32,09.11.202532
130,45970.33333333334,369.681,0.0,33,09.11.202533
131,45970.34375,407.384,0.0,34,09.11.202534
132,45970.35416666666,445.87,0.0,35,09.11.202535
133,45970.36458333334,483.062,0.0,36,09.11.202536
134,45970.375,585.105,0.0,37,09.11.202537
135,45970.38541666666,622.101,0.0,38,09.11.202538
136,45970.39583333334,656.568,0.0,39,09.11.202539
137,45970.40625,689.69,0.0,40,09.11.202540
138,45970.41666666666,735.692,0.0,41,09.11.202541
139,45970.42708333334,758.768,0.0,42,09.11.202542
140,45970.4375,775.012,0.0,43,09.11.202543
141,45970.44791666666,787.659,0.0,44,09.11.202544
142,45970.45833333334,793.347,0.0,45,09.11.202545
143,45970.46875,792.462,0.0,46,09.11.202546
144,45970.47916666666,784.09,0.0,47,09.11.202547
145,45970.48958333334,769.738,0.0,48,09.11.202548
146,45970.5,734.493,0.0,49,09.11.202549
147,45970.51041666666,709.548,0.0,50,09.11.202550
148,45970.52083333334,679.659,0.0,51,09.11.202551
149,45970.53125,648.005,0.0,52,09.11.202552
150,45970.54166666666,574.587,0.0,53,09.11.202553
151,45970.55208333334,535.737,0.0,54,09.11.202554
152,45970.5625,495.796,0.0,55,09.11.202555
153,45970.57291666666,455.69,0.0,56,09.11.202556
154,45970.58333333334,337.84,0.0,57,09.11.202557
155,45970.59375,290.875,0.0,58,09.11.202558
156,45970.60416666666,247.358,0.0,59,09.11.202559
157,45970.61458333334,208.867,0.0,60,09.11.202560
158,45970.625,97.677,0.0,61,09.11.202561
159,45970.63541666666,67.313,0.0,62,09.11.202562
160,45970.64583333334,45.365,0.0,63,09.11.202563
161,45970.65625,29.939,0.0,64,09.11.202564
162,45970.66666666666,9.429,0.0,65,09.11.202565
163,45970.67708333334,7.498,0.0,66,09.11.202566
164,45970.6875,7.557,0.0,67,09.11.202567
165,45970.69791666666,7.277,0.0,68,09.11.202568
166,45970.70833333334,7.606,0.0,69,09.11.202569
167,45970.71875,7.79,0.0,70,09.11.202570
168,45970.72916666666,7.61,0.0,71,09.11.202571
169,45970.73958333334,7.574,0.0,72,09.11.202572
170,45970.75,7.046,0.0,73,09.11.202573
171,45970.76041666666,0.0,0.0,74,09.11.202574
172,45970.77083333334,6.656,0.0,75,09.11.202575
173,45970.78125,4.156,0.0,76,09.11.202576
174,45970.79166666666,0.58,0.0,77,09.11.202577
175,45970.80208333334,0.0,0.0,78,09.11.202578
176,45970.8125,0.0,0.0,79,09.11.202579
177,45970.82291666666,0.0,0.0,80,09.11.202580
178,45970.83333333334,0.0,0.0,81,09.11.202581
179,45970.84375,0.0,0.0,82,09.11.202582
180,45970.85416666666,0.0,0.0,83,09.11.202583
181,45970.86458333334,0.0,0.0,84,09.11.202584
182,45970.875,0.54,0.0,85,09.11.202585
183,45970.88541666666,0.0,0.0,86,09.11.202586
184,45970.89583333334,0.0,0.0,87,09.11.202587
185,45970.90625,0.0,0.0,88,09.11.202588
186,45970.91666666666,0.0,0.0,89,09.11.202589
187,45970.92708333334,0.0,0.0,90,09.11.202590
188,45970.9375,0.0,0.0,91,09.11.202591
189,45970.94791666666,0.0,0.0,92,09.11.202592
190,45970.95833333334,0.0,0.0,93,09.11.202593
191,45970.96875,0.0,0.0,94,09.11.202594
192,45970.97916666666,0.0,0.0,95,09.11.202595
193,45970.98958333334,0.0,0.0,96,09.11.202596
"@

$lines = $csvData -split "`n"
$count = 0
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line.Split(",")
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 1).Value = [double]$parts[1]
    $ws.Cells.Item($r, 2).Value = [double]$parts[2]
    $ws.Cells.Item($r, 3).Value = [double]$parts[3]
    $ws.Cells.Item($r, 4).Value = [int]$parts[4]
    $ws.Cells.Item($r, 5).Value = [string]$parts[5]
    $count = $count + 1
}

Write-Host "Updated $count rows"
